$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the local path of the citizenship document's image next to its
# metadata row (column P = "image").
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Citizenship/citizenship.jpg"
